# Update existing "Temps passé" (Louis) values that were still marked as in-progress
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E6").Value = "10h (toujour en cours)"
$ws.Range("E8").Value = "2h "
$ws.Range("E9").Value = "8h"

# New task row (row 11) for Tommy: "magnétomètre ", mirroring the formatting
# used by the existing data rows (A:E block fill).
$ws.Range("A10:E10").Copy()
$ws.Range("A11:E11").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("A11").Value = "magnétomètre "
$ws.Range("C11").Value = "2h"
$ws.Range("E11").Value = "2h"
